$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Update the "FilesTab" Cypher query text in B4: remove the "File Type" and
# "Breed" coalesce lines from the RETURN clause.
$newQuery = @"
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['French Bulldog']
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS ``File Name``,
        coalesce(labels(parent)[0], '') AS ``Association``,
        coalesce(f.file_description, '') AS ``Description``,
        coalesce(f.file_format, '') AS ``Format``,
        coalesce(f.file_size, '') AS ``Size``,
        coalesce(c.case_id, '') AS ``Case ID``,
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS ``Study Code``
"@

$ws.Range("B4").Value = $newQuery

# Row 4 shrinks now that it has two fewer wrapped lines of text.
$ws.Rows.Item(4).RowHeight = 217.5

# Selection moves to B4 (the cell that was just edited).
$ws.Activate()
$ws.Range("B4").Select()
